$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price / Volume(1h) columns (and the two swapped Stacks/Bittensor
# rows) from the latest coinranking.com data pull.
#
# All of these cells hold plain TEXT in the workbook (e.g. "60.944.26" or
# "  +1.01%  "), not numbers. A few of the new Price values (like "7.40" or
# "0.0000144") look like ordinary numbers/decimals though, and Excel would
# silently convert a bare `.Value = "7.40"` into the number 7.4. To keep those
# cells as text (matching the source data), this writes them with a leading
# apostrophe - Excel's normal "force text" entry marker - then clears the
# quote-prefix formatting that gesture applies so only the cell text changes.

$ws.Range("D2").Value = '60.944.26'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.637.98'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("D5").Value = '''529.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = '''155.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").Value = '''6.66'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("E10").Value = '  +5.05%  '
$ws.Range("D11").Value = '''0.351'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '3.098.40'
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = '60.929.49'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '''21.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.12%  '
$ws.Range("D16").Value = '''0.0000144'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("D17").Value = '2.645.09'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").Value = '''4.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").Value = '''353.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '''10.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").Value = '''6.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.76%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '''61.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.06%  '
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("D27").Value = '0.0₃0864'
$ws.Range("E27").Value = '  +3.65%  '
$ws.Range("D28").Value = '''7.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '''6.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.65%  '
$ws.Range("D31").Value = '''19.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.79%  '
$ws.Range("E32").Value = '  +4.06%  '
$ws.Range("D33").Value = '''150.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").Value = '''4.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.88%  '
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").Value = '''0.927'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.50%  '
$ws.Range("D37").Value = '''0.893'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.21%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '''308.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.86%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''1.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("D40").Value = '''3.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("D42").Value = '''0.641'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.26%  '
$ws.Range("E43").Value = '  +1.89%  '
$ws.Range("D44").Value = '''0.0563'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '''19.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("D49").Value = '''19.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.40%  '
$ws.Range("D51").Value = '1.985.12'
$ws.Range("E51").Value = '  -0.18%  '
